# Updates vm_pu.xlsx res_bus results for the 380 kV case (B changed from 1.05 -> 1.02
# slack/reference voltage, with the dependent bus voltages recomputed accordingly).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","I","J","K","L","M","N")

# Each inner array holds the new values for columns B,C,D,E,F,I,J,K,L,M,N
# (column G stays 1, column H is blank) for data rows 2..25.
$data = @(
    @(1.02,1.070012994947984,1.070204931461848,1.07372751879298,1.079916712169917,1.055512826316385,1.074944804679304,1.072905391148903,1.076418612711777,1.082591510264475,1.076471349599841),
    @(1.02,1.071358866210923,1.071261060445785,1.075025283657365,1.081116722468656,1.055956795270299,1.075946184319272,1.073777243667104,1.077532185862855,1.083608752920098,1.077474151313758),
    @(1.02,1.072229102381869,1.071943822924284,1.075864715919043,1.08189287460222,1.056242509706581,1.076593001839951,1.074340163395683,1.078251878088413,1.084266064734413,1.078121887389507),
    @(1.02,1.072594802052421,1.072230709196117,1.076217542528426,1.082219091353511,1.056362251086269,1.076864653123126,1.074576523221719,1.078554232655951,1.08454218259889,1.078393924448663),
    @(1.02,1.072656196090709,1.072278870103283,1.076276779612365,1.082273860104313,1.056382334356131,1.076910248764649,1.074616192037837,1.078604987434702,1.084588531404168,1.078439584841227),
    @(1.02,1.07223398945404,1.071947656887836,1.075869430679441,1.081897233828193,1.056244111160711,1.076596632718515,1.074343322790469,1.078255918963194,1.084269755080717,1.078125523424335),
    @(1.02,1.070467971710471,1.070561985652565,1.074166170235179,1.080322332116599,1.055663192123225,1.075283463663085,1.073200293177694,1.076795130579119,1.082935482484382,1.076810489518223),
    @(1.02,1.067351007323821,1.06811538015372,1.071162305990664,1.077544475917912,1.054627510123995,1.072960627857934,1.071176633580329,1.074214290650522,1.080577231103883,1.074484355019796),
    @(1.02,1.065269417948613,1.066480889875352,1.069157850084183,1.075690609083973,1.053928888044947,1.071405944976645,1.069821005649345,1.072489017482178,1.079000144260134,1.072927464310492),
    @(1.02,1.064367155587211,1.065772294106724,1.068289408025262,1.074887361954635,1.053624421532582,1.070731259664942,1.069232425793059,1.071740796812601,1.078316047286397,1.072251820868259),
    @(1.02,1.064031872134478,1.065508959407988,1.067966751006077,1.074588920236856,1.053511033227765,1.070480423215025,1.069013560273871,1.071462695589558,1.078061758692025,1.072000628201811),
    @(1.02,1.064103798122589,1.065565451554626,1.068035965669019,1.074652940669774,1.053535368819867,1.070534238868898,1.069060518609998,1.071522357392383,1.0781163128575,1.072054520280085),
    @(1.02,1.064339443897927,1.065750529460883,1.06826273872923,1.074862694319431,1.053615054865406,1.070710530121954,1.069214339228364,1.071717812545956,1.078295031493602,1.072231061886941),
    @(1.02,1.064484613940852,1.065864544664692,1.068402450513531,1.074991919771097,1.053664112811216,1.07081911863171,1.069309081202625,1.071838215123395,1.078405121402149,1.072339804604837),
    @(1.02,1.065329278332285,1.066527898874815,1.069215474865451,1.075743906915449,1.053949053080842,1.071450689787741,1.06986003414613,1.07253864955286,1.079045519814054,1.072972272664351),
    @(1.02,1.065858864213842,1.066943773666677,1.06972532707626,1.076215469593593,1.054127263072889,1.071846454459444,1.070205206331106,1.072977698875197,1.079446899273365,1.073368599367283),
    @(1.02,1.06616767421408,1.067186264513125,1.070022667325604,1.076490475086611,1.054231021086718,1.072077153048735,1.07040638678099,1.073233676657563,1.079680900777013,1.073599625575032),
    @(1.02,1.066272955578422,1.067268933822648,1.070124044657186,1.076584236633163,1.054266367903452,1.072155790938877,1.070474958307426,1.073320939532642,1.079760669599872,1.073678375139999),
    @(1.02,1.065802053848959,1.066899162728811,1.069670629745586,1.07616488047822,1.054108162374382,1.071804007601075,1.070168188432279,1.072930604676611,1.079403847113105,1.073326092229506),
    @(1.02,1.064270056077253,1.065696032260996,1.068195961926875,1.07480092933142,1.053591597492193,1.07065862307041,1.069169049558463,1.071660260838023,1.078242408427796,1.07217908112143),
    @(1.02,1.063305998328398,1.064938818875674,1.067268320874315,1.073942894804303,1.053265099929133,1.069937151914384,1.06853945698131,1.070860509263005,1.077511096745557,1.071456585393609),
    @(1.02,1.063817143842568,1.065340304783743,1.067760125761708,1.074397800245769,1.053438345358792,1.070319743889594,1.068873349012165,1.071284572100446,1.077898881258733,1.071839720693308),
    @(1.02,1.065827724280653,1.066919320744315,1.069695345262818,1.076187739673312,1.054116793741346,1.07182318795484,1.070184915698898,1.072951884864532,1.079423300890558,1.073345299821574),
    @(1.02,1.068157436449925,1.068748478670831,1.071939193673523,1.078262952342972,1.054896692410985,1.073562204078385,1.07170093744879,1.074882317216495,1.081187751582958,1.075086785547483)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $addr = $cols[$j] + $row
        $ws.Range($addr).Value = $vals[$j]
    }
}
